$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 9) with a new pair of addresses/coordinates.
# The coordinate values are stored as literal text (the quote characters
# are part of the text itself), so they are not numeric and Excel keeps
# them as plain text without needing any special formatting.
$ws.Range("A9").Value = "Av. 13 740"
$ws.Range("B9").Value = '"-34.918351696395646"'
$ws.Range("C9").Value = '"-57.958368194341084"'
$ws.Range("D9").Value = "Av. 13 716"
$ws.Range("E9").Value = '"-34.91798918721827"'
$ws.Range("F9").Value = '"-57.9588346727334"'

# Update the selected cell to match the committed state
$ws.Range("F11").Select()
